$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing "0x4000 / RTD Error Code"
# row (and everything below it) down by one.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new message ID entry.
$ws.Range("A2").Value = "0x4001"
$ws.Range("B2").Value = "Main Average Update Time"

# Restore the active selection to match the saved file.
$ws.Range("A2").Select()
